$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook lists DCPS schools that operate "swing space" co-locations,
# e.g. "Bancroft Elementary School @ Sharpe". Strip the " @ <host>" /
# "(...)" suffix so the roster shows the school's own name. Shared strings
# that become unused are dropped by the engine and the new (shorter)
# names are appended to the shared-string table, which also re-homes each
# edited row onto the new string index.
$renames = @{
    5  = "Bancroft Elementary School"
    10 = "Bruce-Monroe Elementary School"
    24 = "Hyde-Addison Elementary School"
    29 = "Kimball Elementary School"
    36 = "Malcolm X Elementary School"
    42 = "Murch Elementary School"
    54 = "School-Within-School"
    69 = "Watkins Elementary School"
}

foreach ($row in ($renames.Keys | Sort-Object)) {
    $ws.Cells.Item($row, 2).Value = $renames[$row]
}

# Restore the active selection to B3:B70 (anchored at B3), matching the
# sheetView saved with this upload.
$ws.Range("B3:B70").Select()
